$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date (column G, row 2)
$wsOverview.Range("G2").Value = "2016-08-14 17:32:50"

# zh-cn sheet: Correspond Handoff Datetime (H2), Correspond Handback DateTime (K2)
$wsZhCn.Range("H2").Value = "2016-08-14 17:32:42"
$wsZhCn.Range("K2").Value = "2016-08-14 17:33:14"

# de-de sheet: Correspond Handoff Datetime (H2)
$wsDeDe.Range("H2").Value = "2016-08-14 17:33:24"
